$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.1844109867082959
    "C2" = 0.2434303829804775
    "D2" = 0.7550547198149107
    "E2" = 0.6915081552015099
    "F2" = 0.4175641082963374
    "G2" = 0.5230975067184298
    "H2" = 0.6743329343542892
    "I2" = 0.188602888553233
    "J2" = 0.4647918628642795
    "K2" = 0.2901314376394514
    "B3" = 0.2658583626134861
    "C3" = 0.7576582593773593
    "D3" = 0.6694404695474946
    "E3" = 0.4115821279919825
    "F3" = 0.5177542565360187
    "G3" = 0.663964400983923
    "H3" = 0.1801255664296092
    "I3" = 0.4570195545574687
    "J3" = 0.2815074641064828
    "K3" = 0.589805676679554
    "B4" = 0.7790828308345794
    "C4" = 0.7896183183467069
    "D4" = 0.3223877835781994
    "E4" = 0.4880272526874408
    "F4" = 0.6725905210138331
    "G4" = 0.154185393004174
    "H4" = 0.4342409915292412
    "I4" = 0.267889428957865
    "J4" = 0.5712848794943781
    "K4" = -0.04166000770222517
    "B5" = 0.7472976860263175
    "C5" = 0.2931833097224661
    "D5" = 0.4850066528825197
    "E5" = 0.6537292374023527
    "F5" = 0.1338542619389697
    "G5" = 0.4197108360949934
    "H5" = 0.251248383465597
    "I5" = 0.553701034065628
    "J5" = -0.05813389123226798
    "K5" = 0.6310880987550094
    "B6" = 0.6333810503403763
    "C6" = 0.5603777146664568
    "D6" = 0.4626359510893536
    "E6" = 0.1585481745257029
    "F6" = 0.4290902396512893
    "G6" = 0.196437584854206
    "H6" = 0.5348069431769582
    "I6" = -0.07299789204731372
    "J6" = 0.6023763881256388
    "K6" = 0.3304428394308724
    "B7" = 1.011072586656707
    "C7" = 0.509746653301712
    "D7" = -0.08173505626510533
    "E7" = 0.4650353207171521
    "F7" = 0.1943868274146131
    "G7" = 0.458238703740928
    "H7" = -0.1003972864225971
    "I7" = 0.5760873379279133
    "J7" = 0.2870455487175398
    "B8" = 0.8220722402252505
    "C8" = 0.051553683470419
    "D8" = 0.2859568989009605
    "E8" = 0.2225810661851803
    "F8" = 0.4942640149333215
    "G8" = -0.1379388566268107
    "H8" = 0.5705449988405521
    "I8" = 0.2912208776562884
    "B9" = 0.2871441745782602
    "C9" = 0.3706356397752701
    "D9" = 0.07698069955146319
    "E9" = 0.5042177577925642
    "F9" = -0.1225736869272658
    "G9" = 0.5316966638831291
    "H9" = 0.2753750686291025
    "B10" = 0.6816598262566529
    "C10" = 0.1940694739626584
    "D10" = 0.3417536163429973
    "E10" = -0.0936147364620642
    "F10" = 0.5674395363380327
    "G10" = 0.2440474222454754
    "B11" = 0.4408946513667728
    "C11" = 0.3594094838808868
    "D11" = -0.1882369755730587
    "E11" = 0.5996569245865127
    "F11" = 0.2564355480731927
    "B12" = 0.5990858432970987
    "C12" = -0.1032127321038452
    "D12" = 0.4831723462284986
    "E12" = 0.2715408197250452
    "B13" = 0.0616473449302421
    "C13" = 0.4967096184764148
    "D13" = 0.2085679007350822
    "B14" = 0.7505586603418228
    "C14" = 0.3078859509171186
    "B15" = 0.3519456421565676
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

$clearCells = @("J8", "I9", "H10", "G11", "F12", "E13", "D14", "C15", "B16")
foreach ($addr in $clearCells) {
    $ws.Range($addr).ClearContents()
}

Write-Host "Done applying changes"
